$d = $word.ActiveDocument

# 1) Ativacao date: 01/01/2012 -> 01/01/2021
$d.Content.Find.Execute('Ativação: 01/01/2012', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2021', 2) | Out-Null

# 2) Objetivos paragraph full replacement
$d.Content.Find.Execute('Habilitar o aluno resolver problemas práticos em mecânica dos fluidos, modelando situações reais (através das equações de conservação e fenomenológicas), promovendo abstrações e adequando os casos ilustrados a novas situações. Capacitar o aluno a realizar cálculos de transferência de calor, utilizando os mecanismos, de condução e convecção, combinados ou não. Fornecer ao aluno condições de estabelecer um pré dimensionamento de um trocador de calor.', $true, $false, $false, $false, $false, $true, 1, $false, '1.Integrar o aluno com o conhecimento teórico e prático de fenômenos de transporte, abordando temas relacionados à Mecânica dos Fluidos, Transferência de Calor e Transferência de Massa, de forma a torná-lo capaz de compreender os processos físicos envolvidos nesses processos;2.Habilitar o aluno resolver problemas práticos em mecânica dos fluidos, modelando situações reais (através das equações de conservação e fenomenológicas), promovendo abstrações e adequando os casos ilustrados a novas situações.3.Capacitar o aluno a realizar cálculos de transferência de calor, utilizando os mecanismos, de condução e convecção, combinados ou não.', 2) | Out-Null

# 3) Programa resumido paragraph full replacement
$d.Content.Find.Execute('Conceitos ligados ao escoamento de fluídos e equações fundamentais, Escoamento incompressível de fluidos não viscosos, Escoamento viscoso incompressível, Transferência de Calor.', $true, $false, $false, $false, $false, $true, 1, $false, 'Sistemas de unidades e dimensões. Conceitos básicos de escoamentos de fluidos e analogias entre transferência de momento linear, de calor e de massa. Estática e cinemática dos fluídos.  Descrição e classificação de escoamentos. Introdução à análise de escoamentos (volume de controle e diferencial). Introdução aos modos de transferência de calor (regime permanente e transiente). Introdução à transferência de massa.', 2) | Out-Null

# 4) Programa section: collapse <w:br/> runs into one run, then replace whole text
$progPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("1. Conceitos ligados ao escoamento")) {
        $progPara = $p
        break
    }
}
$progPara.Range.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$progPara.Range.Text = '1. Conceitos ligados ao escoamento de fluídos e equações fundamentais 1.1. Características e definições dos escoamentos;1.2. Conceitos de sistema e volume de controle;1.3. Equação da conservação da massa;1.4. Equação da conservação da energia;1.5. Equação da conservação da quantidade de movimento; 1.6. Introdução à análise diferencial do movimento de fluidos.2. Escoamento incompressível de fluidos não viscosos2.1. Equação de Euler;2.2. Equação de Bernoulli;2.3. Aplicações da equação de Bernoulli.3. Escoamento viscoso incompressível3.1. Atrito e perda de carga;3.2. Avaliação das perdas de carga: regime laminar e turbulento;3.3. Equações de Hagen - Poiseuille e Darcy  Weisbach3.4. Diagrama de Moody e Moody Rouse;3.5. Método dos comprimentos equivalentes.3.6. Presença de máquina no escoamento (bomba e turbina), Potência e rendimento;3.7. Medidores de vazão.4. Transferência de Calor4.1. Definição de Calor.4.2. Mecanismo da Condução.4.3. Mecanismo da Convecção.4.4. Associação de Mecanismos.5. Transferência de Massa5.1. Difusão e convecção mássica;5.2. 1ª lei de Fick;5.3. Concentrações mássica e molar;5.4. Frações mássica e molar;5.5. Velocidades médias mássica e molar;5.6. Fluxos difusivo mássico, difusivo molar, convectivo mássico e convectivo molar;5.7. Fluxo mássico total e fluxo molar total.'

# 5) Metodo run replacement
$d.Content.Find.Execute('Aulas expositivas teóricas, aulas práticas, aulas de exercícios.', $true, $false, $false, $false, $false, $true, 1, $false, 'Duas provas escritas', 2) | Out-Null

# 6) Bibliografia section: collapse <w:br/> runs into one run, then replace whole text
$bibPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("1. FOX, R.W.")) {
        $bibPara = $p
        break
    }
}
$bibPara.Range.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$bibPara.Range.Text = '1.LEVI, Celso P., Fundamentos de Fenômenos de Transporte. 2 ed. Ed. Gen/LTC. Rio de Janeiro/RJ. 2013.2.FOX, Robert.W., McDonald, Alan.T.; PRITCHARD, Philip J., Introdução à Mecânica dos Fluidos, 7 ed. Ed. Gen/LTC. Rio de Janeiro/RJ. 2012.3.KREITH, Frank; BOHN, Mark S., Fundamentos de Transferência de Calor, Ed. Cengage Learning, São Paulo/SP. 2011. 4.INCROPERA/ DEWITT/ BERGMAN/ LAVINE. Fundamentos de Transferência de Calor e Massa. 7 ed. Ed. Gen/LTC. Rio de Janeiro/RJ. 20145.CREMASCO, Marco. A., Fundamentos de Transferência de Massa, 3ed. Ed.  Blücher, São Paulo/SP. 2016'
